$wb = $excel.ActiveWorkbook

# Sheet "PFOS-Tia" (2nd tab) - new "PFOA-Tia" sheet gets inserted right after it
$afterSheet = $wb.Worksheets.Item("PFOS-Tia")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "PFOA-Tia"

$newSheet.Range("A1").Value = "Replicates"
$newSheet.Range("B1").Value = "Genotypes"
$newSheet.Range("C1").Value = "Treatment"
$newSheet.Range("D1").Value = "Age_maturity"
$newSheet.Range("E1").Value = "Day_1brood"
$newSheet.Range("F1").Value = "first_brood"
$newSheet.Range("G1").Value = "Size_maturity"
$newSheet.Range("H1").Value = "Day_2brood"
$newSheet.Range("I1").Value = "Second_brood"
$newSheet.Range("J1").Value = "Fecundity"
$newSheet.Range("K1").Value = "Invterval_brood"
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "LRV-0-1"
$newSheet.Range("C2").Value = "control"
$newSheet.Range("D2").Value = 8
$newSheet.Range("E2").Value = 10
$newSheet.Range("F2").Value = 26
$newSheet.Range("G2").Value = 2590.66
$newSheet.Range("H2").Value = 13
$newSheet.Range("I2").Value = 26
$newSheet.Range("J2").Formula = "=(F2+I2)"
$newSheet.Range("K2").Value = 3
$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = "LRV-0-1"
$newSheet.Range("C3").Value = "control"
$newSheet.Range("D3").Value = 7
$newSheet.Range("E3").Value = 10
$newSheet.Range("F3").Value = 24
$newSheet.Range("G3").Value = 2583.37
$newSheet.Range("H3").Value = 12
$newSheet.Range("I3").Value = 20
$newSheet.Range("J3").Formula = "=(F3+I3)"
$newSheet.Range("K3").Value = 2
$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = "LRV-0-1"
$newSheet.Range("C4").Value = "control"
$newSheet.Range("D4").Value = 8
$newSheet.Range("E4").Value = 11
$newSheet.Range("F4").Value = 26
$newSheet.Range("G4").Value = 2601.63
$newSheet.Range("H4").Value = 13
$newSheet.Range("I4").Value = 26
$newSheet.Range("J4").Formula = "=(F4+I4)"
$newSheet.Range("K4").Value = 2
$newSheet.Range("A5").Value = 4
$newSheet.Range("B5").Value = "LRV-0-1"
$newSheet.Range("C5").Value = "control"
$newSheet.Range("D5").Value = 8
$newSheet.Range("E5").Value = 10
$newSheet.Range("F5").Value = 28
$newSheet.Range("G5").Value = 2650.88
$newSheet.Range("H5").Value = 12
$newSheet.Range("I5").Value = 19
$newSheet.Range("J5").Formula = "=(F5+I5)"
$newSheet.Range("K5").Value = 2
$newSheet.Range("A6").Value = 1
$newSheet.Range("B6").Value = "LRV-0-1"
$newSheet.Range("C6").Value = "PFOA"
$newSheet.Range("D6").Value = 7
$newSheet.Range("E6").Value = 12
$newSheet.Range("F6").Value = 9
$newSheet.Range("G6").Value = 2310.09
$newSheet.Range("H6").Value = 15
$newSheet.Range("I6").Value = 9
$newSheet.Range("J6").Formula = "=(F6+I6)"
$newSheet.Range("K6").Value = 3
$newSheet.Range("A7").Value = 2
$newSheet.Range("B7").Value = "LRV-0-1"
$newSheet.Range("C7").Value = "PFOA"
$newSheet.Range("D7").Value = 8
$newSheet.Range("E7").Value = 13
$newSheet.Range("F7").Value = 8
$newSheet.Range("G7").Value = 2290.9899999999998
$newSheet.Range("H7").Value = 17
$newSheet.Range("I7").Value = 8
$newSheet.Range("J7").Formula = "=(F7+I7)"
$newSheet.Range("K7").Value = 4
$newSheet.Range("A8").Value = 3
$newSheet.Range("B8").Value = "LRV-0-1"
$newSheet.Range("C8").Value = "PFOA"
$newSheet.Range("D8").Value = 8
$newSheet.Range("E8").Value = 11
$newSheet.Range("F8").Value = 7
$newSheet.Range("G8").Value = 2295.6999999999998
$newSheet.Range("H8").Value = 15
$newSheet.Range("I8").Value = 10
$newSheet.Range("J8").Formula = "=(F8+I8)"
$newSheet.Range("K8").Value = 4
$newSheet.Range("A9").Value = 4
$newSheet.Range("B9").Value = "LRV-0-1"
$newSheet.Range("C9").Value = "PFOA"
$newSheet.Range("D9").Value = 9
$newSheet.Range("E9").Value = 13
$newSheet.Range("F9").Value = 12
$newSheet.Range("G9").Value = 2203.0300000000002
$newSheet.Range("H9").Value = 17
$newSheet.Range("I9").Value = 8
$newSheet.Range("J9").Formula = "=(F9+I9)"
$newSheet.Range("K9").Value = 4
$newSheet.Range("A10").Value = 1
$newSheet.Range("B10").Value = "LR2-36-01"
$newSheet.Range("C10").Value = "control"
$newSheet.Range("D10").Value = 6
$newSheet.Range("E10").Value = 9
$newSheet.Range("F10").Value = 28
$newSheet.Range("G10").Value = 2539.61
$newSheet.Range("H10").Value = 12
$newSheet.Range("I10").Value = 19
$newSheet.Range("J10").Formula = "=(F10+I10)"
$newSheet.Range("K10").Value = 3
$newSheet.Range("A11").Value = 2
$newSheet.Range("B11").Value = "LR2-36-01"
$newSheet.Range("C11").Value = "control"
$newSheet.Range("D11").Value = 6
$newSheet.Range("E11").Value = 9
$newSheet.Range("F11").Value = 26
$newSheet.Range("G11").Value = 2544.69
$newSheet.Range("H11").Value = 12
$newSheet.Range("I11").Value = 28
$newSheet.Range("J11").Formula = "=(F11+I11)"
$newSheet.Range("K11").Value = 3
$newSheet.Range("A12").Value = 3
$newSheet.Range("B12").Value = "LR2-36-01"
$newSheet.Range("C12").Value = "control"
$newSheet.Range("D12").Value = 8
$newSheet.Range("E12").Value = 10
$newSheet.Range("F12").Value = 29
$newSheet.Range("G12").Value = 2606.88
$newSheet.Range("H12").Value = 13
$newSheet.Range("I12").Value = 22
$newSheet.Range("J12").Formula = "=(F12+I12)"
$newSheet.Range("K12").Value = 3
$newSheet.Range("A13").Value = 4
$newSheet.Range("B13").Value = "LR2-36-01"
$newSheet.Range("C13").Value = "control"
$newSheet.Range("D13").Value = 6
$newSheet.Range("E13").Value = 9
$newSheet.Range("F13").Value = 19
$newSheet.Range("G13").Value = 2568.98
$newSheet.Range("H13").Value = 12
$newSheet.Range("I13").Value = 28
$newSheet.Range("J13").Formula = "=(F13+I13)"
$newSheet.Range("K13").Value = 3
$newSheet.Range("A14").Value = 1
$newSheet.Range("B14").Value = "LR2-36-01"
$newSheet.Range("C14").Value = "PFOA"
$newSheet.Range("D14").Value = 8
$newSheet.Range("E14").Value = 12
$newSheet.Range("F14").Value = 9
$newSheet.Range("G14").Value = 2271.46
$newSheet.Range("H14").Value = 17
$newSheet.Range("I14").Value = 20
$newSheet.Range("J14").Formula = "=(F14+I14)"
$newSheet.Range("K14").Value = 5
$newSheet.Range("A16").Value = 3
$newSheet.Range("B16").Value = "LR2-36-01"
$newSheet.Range("C16").Value = "PFOA"
$newSheet.Range("D16").Value = 8
$newSheet.Range("E16").Value = 11
$newSheet.Range("F16").Value = 8
$newSheet.Range("G16").Value = 2290.2800000000002
$newSheet.Range("H16").Value = 17
$newSheet.Range("I16").Value = 9
$newSheet.Range("J16").Formula = "=(F16+I16)"
$newSheet.Range("K16").Value = 6
$newSheet.Range("A17").Value = 4
$newSheet.Range("B17").Value = "LR2-36-01"
$newSheet.Range("C17").Value = "PFOA"
$newSheet.Range("D17").Value = 8
$newSheet.Range("E17").Value = 13
$newSheet.Range("F17").Value = 9
$newSheet.Range("G17").Value = 2201.0100000000002
$newSheet.Range("H17").Value = 17
$newSheet.Range("I17").Value = 17
$newSheet.Range("J17").Formula = "=(F17+I17)"
$newSheet.Range("K17").Value = 4

# Row 15 - special row: highlighted yellow, only A-D populated, E:J empty but styled, no K value
$newSheet.Range("A15:J15").Interior.Color = 65535
$newSheet.Range("A15").Value = 2
$newSheet.Range("B15").Value = "LR2-36-01"
$newSheet.Range("C15").Value = "PFOA"
$newSheet.Range("D15").Value = 8

# Selection state to match target
$newSheet.Range("B17").Select()
